$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("137:137").Insert()

$ws.Range("A137").Value = 10
$ws.Range("B137").Value = "Vega Modelo de Temuco"
$ws.Range("C137").Value = "La Araucanía"
$ws.Range("D137").Value = 44977
$ws.Range("E137").Value = 9
$ws.Range("F137").Value = 100114007
$ws.Range("G137").Value = "Jengibre"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 25
$ws.Range("K137").Value = 35000
$ws.Range("L137").Value = 35000
$ws.Range("M137").Value = 35000
$ws.Range("N137").Value = "$/caja 13 kilos"
$ws.Range("O137").Value = "Perú"
$ws.Range("P137").Value = 2692
$ws.Range("Q137").Value = 13
$ws.Range("R137").Value = "Hortaliza"
